$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying external/query data was refreshed and the table was
# re-sorted ascending by "Name" (natural numeric order: bin1, bin2, ... bin17,
# bins, frequency, id, max, mean, min, size). Re-write A2:B25 in that order.
# (bin17's value also changed from 0 -> 1 as part of the refreshed data.)
$names  = @("bin1","bin2","bin3","bin4","bin5","bin6","bin7","bin8","bin9","bin10","bin11","bin12","bin13","bin14","bin15","bin16","bin17","bins","frequency","id","max","mean","min","size")
$values = @(95,77,33,31,16,16,8,6,6,4,3,0,1,2,0,1,1,17,1.7275294117647058,3,29.375,4.6044166666666664,0.007,300)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Input modeling should start -> add a running total next to the bin
# frequencies so downstream model translation can reference it.
$ws.Range("C1").Formula = "=SUM(B2:B18)"

# Leave the selection on the newly added cell.
$ws.Range("C1").Select() | Out-Null
